$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 47 and 48: Aptos and Maker swap positions with updated price/volume
$ws.Range("B47").Value = "Maker"
$ws.Range("C47").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D47").Value = "1.008.03"
$ws.Range("E47").Value = "  +5.10%  "

$ws.Range("B48").Value = "Aptos"
$ws.Range("C48").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D48").Value = "'7.578"
$ws.Range("E48").Value = "  -0.90%  "

$ws.Range("D2").Value = "29.984.15"
$ws.Range("E2").Value = "  -0.10%  "
$ws.Range("D3").Value = "1.896.49"
$ws.Range("E3").Value = "  -0.71%  "
$ws.Range("D4").Value = "'0.9999"
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").Value = "'0.8397"
$ws.Range("E5").Value = "  +5.16%  "
$ws.Range("D6").Value = "'241.61"
$ws.Range("E6").Value = "  -0.23%  "
$ws.Range("D7").Value = "'0.9996"
$ws.Range("E7").Value = "  -0.04%  "
$ws.Range("D8").Value = "'0.3288"
$ws.Range("E8").Value = "  +3.72%  "
$ws.Range("D9").Value = "'26.63"
$ws.Range("E9").Value = "  +1.69%  "
$ws.Range("D10").Value = "'0.07055"
$ws.Range("E10").Value = "  +1.64%  "
$ws.Range("D11").Value = "'0.08073"
$ws.Range("E11").Value = "  +0.45%  "
$ws.Range("D12").Value = "'0.7582"
$ws.Range("E12").Value = "  +1.17%  "
$ws.Range("D13").Value = "1.894.30"
$ws.Range("E13").Value = "  -0.57%  "
$ws.Range("D14").Value = "'5.260"
$ws.Range("E14").Value = "  +0.55%  "
$ws.Range("D15").Value = "'92.26"
$ws.Range("E15").Value = "  -1.13%  "
$ws.Range("D16").Value = "29.972.33"
$ws.Range("E16").Value = "  -0.13%  "
$ws.Range("D17").Value = "'14.10"
$ws.Range("E17").Value = "  +0.64%  "
$ws.Range("D18").Value = "'5.881"
$ws.Range("E18").Value = "  -1.45%  "
$ws.Range("D19").Value = "'244.66"
$ws.Range("E19").Value = "  -2.11%  "
$ws.Range("D20").Value = "'0.000007774"
$ws.Range("E20").Value = "  -0.14%  "
$ws.Range("D21").Value = "'0.9981"
$ws.Range("E21").Value = "  -0.15%  "
$ws.Range("D22").Value = "2.145.60"
$ws.Range("E22").Value = "  +0.07%  "
$ws.Range("D23").Value = "'1.000"
$ws.Range("E23").Value = "  +0.00%  "
$ws.Range("D24").Value = "'6.991"
$ws.Range("E24").Value = "  +0.26%  "
$ws.Range("D25").Value = "'0.1738"
$ws.Range("E25").Value = "  +27.75%  "
$ws.Range("D26").Value = "'9.251"
$ws.Range("E26").Value = "  -0.68%  "
$ws.Range("D27").Value = "'165.85"
$ws.Range("E27").Value = "  -1.68%  "
$ws.Range("D28").Value = "'18.89"
$ws.Range("E28").Value = "  -0.80%  "
$ws.Range("D29").Value = "'2.105"
$ws.Range("E29").Value = "  +2.21%  "
$ws.Range("D30").Value = "'1.359"
$ws.Range("E30").Value = "  -2.16%  "
$ws.Range("D31").Value = "'1.519"
$ws.Range("E31").Value = "  -0.21%  "
$ws.Range("D32").Value = "'0.05915"
$ws.Range("E32").Value = "  +9.82%  "
$ws.Range("D33").Value = "'4.298"
$ws.Range("E33").Value = "  -1.31%  "
$ws.Range("D34").Value = "'4.079"
$ws.Range("E34").Value = "  -1.13%  "
$ws.Range("D35").Value = "'1.275"
$ws.Range("E35").Value = "  +0.96%  "
$ws.Range("D36").Value = "'0.7325"
$ws.Range("E36").Value = "  -0.89%  "
$ws.Range("D37").Value = "'2.721"
$ws.Range("E37").Value = "  -0.18%  "
$ws.Range("D38").Value = "'0.01919"
$ws.Range("E38").Value = "  -0.53%  "
$ws.Range("D39").Value = "'2.773"
$ws.Range("E39").Value = "  -0.50%  "
$ws.Range("D40").Value = "'0.4447"
$ws.Range("E40").Value = "  -0.48%  "
$ws.Range("D41").Value = "'72.42"
$ws.Range("E41").Value = "  -0.86%  "
$ws.Range("D42").Value = "'5.877"
$ws.Range("E42").Value = "  -4.83%  "
$ws.Range("D43").Value = "'0.8427"
$ws.Range("E43").Value = "  +1.37%  "
$ws.Range("D44").Value = "'0.9991"
$ws.Range("E44").Value = "  -0.11%  "
$ws.Range("D45").Value = "'1.886"
$ws.Range("E45").Value = "  -1.28%  "
$ws.Range("D46").Value = "'101.86"
$ws.Range("E46").Value = "  +1.20%  "
$ws.Range("D49").Value = "'9.809"
$ws.Range("E49").Value = "  -0.66%  "
$ws.Range("D50").Value = "2.044.47"
$ws.Range("E50").Value = "  -0.41%  "
$ws.Range("D51").Value = "'35.98"
$ws.Range("E51").Value = "  -1.15%  "
